# Se termina el modulo de pagos
# Registrar el consumo de horas (Dia 8) para las tareas del modulo de pagos
# en las filas 7 y 8 de la hoja "Casos de Uso".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

$ws.Range("AC7").Value = 1
$ws.Range("AC8").Value = 1

# Actualiza la celda activa seleccionada en el panel inferior derecho
$ws.Range("R7").Select() | Out-Null
